$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp banner (A1)
$ws.Range("A1").Value = "Datos actualizados a 6 de Septiembre de 2020 a las 00:29"

# Row 4: Estados Unidos
$ws.Range("B4").Value = 6427291
$ws.Range("C4").Value = 38234
$ws.Range("D4").Value = 3696357
$ws.Range("E4").Value = 2538150
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 673
$ws.Range("H4").Value = 192784

# Row 6: India
$ws.Range("B6").Value = 4110839
$ws.Range("C6").Value = 90600
$ws.Range("D6").Value = 3177673
$ws.Range("E6").Value = 862487
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 1044
$ws.Range("H6").Value = 70679

# Row 9: Colombia
$ws.Range("B9").Value = 658456
$ws.Range("C9").Value = 8394
$ws.Range("D9").Value = 507770
$ws.Range("E9").Value = 129530
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 268
$ws.Range("H9").Value = 21156

# Row 34: Egipto
$ws.Range("B34").Value = 99712
$ws.Range("C34").Value = 130
$ws.Range("D34").Value = 77208
$ws.Range("E34").Value = 16993
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 16
$ws.Range("H34").Value = 5511

# Row 43: Guatemala
$ws.Range("B43").Value = 77481
$ws.Range("C43").Value = 441
$ws.Range("D43").Value = 65595
$ws.Range("E43").Value = 9041
$ws.Range("F43").Value = 0
$ws.Range("G43").Value = 20
$ws.Range("H43").Value = 2845

# Row 47: Japon
$ws.Range("A47").Value = "Japon"
$ws.Range("B47").Value = 70876
$ws.Range("C47").Value = 608
$ws.Range("D47").Value = 61445
$ws.Range("E47").Value = 8082
$ws.Range("F47").Value = 0
$ws.Range("G47").Value = 19
$ws.Range("H47").Value = 1349

# Row 48: Polonia
$ws.Range("A48").Value = "Polonia"
$ws.Range("B48").Value = 70387
$ws.Range("C48").Value = 567
$ws.Range("D48").Value = 52346
$ws.Range("E48").Value = 15928
$ws.Range("F48").Value = 0
$ws.Range("G48").Value = 13
$ws.Range("H48").Value = 2113

# Row 54: Barein
$ws.Range("A54").Value = "Barein"
$ws.Range("B54").Value = 54771
$ws.Range("C54").Value = 676
$ws.Range("D54").Value = 50645
$ws.Range("E54").Value = 3930
$ws.Range("F54").Value = 0
$ws.Range("G54").Value = 1
$ws.Range("H54").Value = 196

# Row 55: Nigeria
$ws.Range("A55").Value = "Nigeria"
$ws.Range("B55").Value = 54743
$ws.Range("C55").Value = 0
$ws.Range("D55").Value = 42816
$ws.Range("E55").Value = 10876
$ws.Range("F55").Value = 0
$ws.Range("G55").Value = 0
$ws.Range("H55").Value = 1051

# Row 108: Malaui
$ws.Range("B108").Value = 5611
$ws.Range("C108").Value = 3
$ws.Range("D108").Value = 3551
$ws.Range("E108").Value = 1885
$ws.Range("F108").Value = 0
$ws.Range("G108").Value = 0
$ws.Range("H108").Value = 175

# Row 113: Suazilandia
$ws.Range("B113").Value = 4819
$ws.Range("C113").Value = 39
$ws.Range("D113").Value = 3903
$ws.Range("E113").Value = 822
$ws.Range("F113").Value = 0
$ws.Range("G113").Value = 0
$ws.Range("H113").Value = 94

# Row 119: Ruanda
$ws.Range("A119").Value = "Ruanda"
$ws.Range("B119").Value = 4349
$ws.Range("C119").Value = 45
$ws.Range("D119").Value = 2199
$ws.Range("E119").Value = 2132
$ws.Range("F119").Value = 0
$ws.Range("G119").Value = 0
$ws.Range("H119").Value = 18

# Row 120: Mozambique
$ws.Range("A120").Value = "Mozambique"
$ws.Range("B120").Value = 4341
$ws.Range("C120").Value = 76
$ws.Range("D120").Value = 2579
$ws.Range("E120").Value = 1736
$ws.Range("F120").Value = 0
$ws.Range("G120").Value = 0
$ws.Range("H120").Value = 26

# Row 154: Togo
$ws.Range("B154").Value = 1477
$ws.Range("C154").Value = 20
$ws.Range("D154").Value = 1094
$ws.Range("E154").Value = 351
$ws.Range("F154").Value = 0
$ws.Range("G154").Value = 1
$ws.Range("H154").Value = 32
